$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.284.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.67%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.920.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.29%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.24%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'373.86"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'103.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.71%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -3.34%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.15%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -5.38%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'37.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.99%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -0.38%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.0839"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.19%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'18.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.91%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'3.377.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.45%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'7.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.76%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'2.913.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.64%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.936"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -8.79%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'51.220.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.92%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'3.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.99%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'7.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.43%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'13.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.65%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.0₃0946"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.15%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'68.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.50%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'261.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.12%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.75%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -5.64%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = "'Dai"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.01%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'25.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.78%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = "'Filecoin"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'7.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.66%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = "'RenderToken"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'6.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +6.28%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "'Hedera"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'0.102"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.80%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'Cosmos"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'9.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.98%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "'Toncoin"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'2.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.41%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "'InjectiveProtocol"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'34.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.66%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = "'OKB"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'51.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.18%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "'FirstDigitalUSD"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.41%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'VeChain"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.0425"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.70%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'LidoDAOToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -6.13%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "'Celestia"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'17.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.46%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'Stacks"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'2.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.10%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'ARBITRUM"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'1.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -6.41%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'Stellar"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.113"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.63%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "'Monero"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'119.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.96%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'22.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.49%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'WEMIXToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.82%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'Maker"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'2.028.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.57%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'ApeXProtocol"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.69%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'NEARProtocol"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'3.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.91%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'TheGraph"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.252"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.11%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'RocketPoolETH"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'3.212.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.24%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'BEAM"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0321"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.84%  "
$ws.Range("E51").Style = "Normal"
